$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. A leading apostrophe forces
# Excel to store the value as text (avoiding numeric auto-conversion),
# and resetting the Style back to "Normal" afterwards strips the resulting
# quote-prefix formatting so the cell keeps the workbook default style.
$updates = [ordered]@{
    "D2" = '''43.299.73'
    "E2" = '''  -1.27%  '
    "D3" = '''2.237.68'
    "E3" = '''  -0.72%  '
    "E4" = '''  +0.01%  '
    "D5" = '''230.18'
    "E5" = '''  +0.16%  '
    "D6" = '''0.640'
    "E6" = '''  +1.59%  '
    "D7" = '''63.97'
    "E7" = '''  +1.05%  '
    "E8" = '''  -0.11%  '
    "D9" = '''0.438'
    "E9" = '''  -0.10%  '
    "D10" = '''0.0950'
    "E10" = '''  -5.68%  '
    "D11" = '''56.35'
    "E11" = '''  -0.07%  '
    "D12" = '''26.59'
    "E12" = '''  +3.92%  '
    "E13" = '''  -1.65%  '
    "D14" = '''2.569.88'
    "E14" = '''  -0.53%  '
    "D15" = '''15.14'
    "E15" = '''  -2.84%  '
    "D16" = '''5.99'
    "E16" = '''  +0.50%  '
    "D17" = '''0.820'
    "E17" = '''  -0.02%  '
    "D18" = '''2.245.46'
    "E18" = '''  -0.68%  '
    "D19" = '''43.167.12'
    "E19" = '''  -1.21%  '
    "D20" = '''0.0₃0959'
    "D21" = '''72.81'
    "E21" = '''  -0.70%  '
    "D22" = '''6.02'
    "E22" = '''  +0.11%  '
    "D23" = '''245.26'
    "E23" = '''  -3.08%  '
    "E24" = '''  +0.06%  '
    "E25" = '''  +26.96%  '
    "D26" = '''2.41'
    "E26" = '''  -0.74%  '
    "E27" = '''  -2.48%  '
    "D28" = '''173.89'
    "E28" = '''  +1.65%  '
    "D29" = '''9.66'
    "E29" = '''  -2.34%  '
    "D30" = '''21.62'
    "E30" = '''  +4.38%  '
    "B31" = '''ImmutableX'
    "C31" = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    "D31" = '''1.40'
    "E31" = '''  +0.97%  '
    "B32" = '''Kaspa'
    "C32" = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    "D32" = '''0.129'
    "E32" = '''  -6.82%  '
    "E33" = '''  +1.57%  '
    "D34" = '''4.88'
    "E34" = '''  +4.20%  '
    "D35" = '''0.0673'
    "E35" = '''  -0.17%  '
    "D36" = '''4.87'
    "E36" = '''  -0.55%  '
    "D37" = '''3.58'
    "E37" = '''  -7.74%  '
    "D38" = '''6.29'
    "E38" = '''  -5.98%  '
    "D39" = '''2.25'
    "E39" = '''  -2.80%  '
    "D40" = '''0.0248'
    "E40" = '''  -0.64%  '
    "E41" = '''  -0.15%  '
    "D42" = '''8.53'
    "E42" = '''  +3.48%  '
    "D43" = '''4.45'
    "E43" = '''  +2.91%  '
    "D44" = '''16.94'
    "E44" = '''  -3.72%  '
    "D45" = '''96.21'
    "E45" = '''  -0.81%  '
    "D46" = '''0.0936'
    "E46" = '''  -2.45%  '
    "E48" = '''  -1.37%  '
    "D49" = '''1.423.46'
    "E49" = '''  -1.68%  '
    "D50" = '''9.74'
    "E50" = '''  +1.71%  '
    "B51" = '''NEARProtocol'
    "C51" = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    "D51" = '''2.26'
    "E51" = '''  -1.70%  '
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
